$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates (table description tweaks) ---
# Etudiants table: the "carte" column is not a primary key, drop the "* " marker
$ws.Range("B4").Value = "v_id_carte"

# Groupes table: first column header should read "* v_id_groupe" (singular) not "* v_id_groupes"
$ws.Range("F13").Value = "* v_id_groupe"

# Creneaux table: the foreign key column now points at "groupe" rather than "classe"
$ws.Range("B7").Value = '$ v_id_groupe'
$ws.Range("B27").Value = '$ v_id_groupe'

# --- View / selection state ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B27").Select()
